$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "30,46 TRY - 60,94 TRY - 609,43 TRY"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("G4").Value = "21,27 TRY - 42,55 TRY - 304,71 TRY"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("G5").Value = "6,09 TRY - 12,19 TRY - 152,35 TRY"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("G6").Value = "4.300,01 TL - 76,17 TL"
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

$ws.Range("G8").Value = "15,23 TRY - 30,47 TRY - 304,71 TRY"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("G9").Value = "10,63 TRY - 21,27 TRY - 152,35 TRY"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("G10").Value = "3,04 TRY - 6,09 TRY - 76,17 TRY"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("G11").Value = "3,04 TRY - 6,09 TRY - 76,17 TRY"
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"

$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"

$ws.Range("G14").Value = "4.300 TL - 6,09 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
